# Update the descr (AlternativeText) attribute on four p:pic shapes so the
# embedded image filenames reflect the new rendered hashes.

$p = $ppt.ActivePresentation

# Slide 3: sales chart picture -> new hash 89f29345
$s = $p.Slides.Item(3)
$s.Shapes.Item(3).AlternativeText = "slide_89f29345_create_sales_chart.png"

# Slide 4: market share picture -> new hash efe173e6
$s = $p.Slides.Item(4)
$s.Shapes.Item(3).AlternativeText = "slide_efe173e6_create_market_share.png"

# Slide 11: growth trend picture -> new hash 27b543ea
$s = $p.Slides.Item(11)
$s.Shapes.Item(3).AlternativeText = "slide_27b543ea_create_growth_trend.png"

# Slide 13: sales chart picture (second occurrence) -> new hash 89f29345
$s = $p.Slides.Item(13)
$s.Shapes.Item(3).AlternativeText = "slide_89f29345_create_sales_chart.png"
